# call of cthulhu product_code
# Adds a new "product_code" column (G) to the checklist sheet, populating
# it with the "No. 128" / "No. 128-N" product codes for the rows that have
# one, widens column F slightly, and leaves the active selection on G12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("G1").Value = "product_code"

# Product codes per row (rows without a code are left untouched). Written
# in this particular order so new shared-string entries land in the same
# sequence as the source edit.
$ws.Range("G4").Value  = "No. 128"
$ws.Range("G15").Value = "No. 128-2"
$ws.Range("G14").Value = "No. 128-8"
$ws.Range("G6").Value  = "No. 128-6"
$ws.Range("G2").Value  = "No. 128-1"
$ws.Range("G13").Value = "No. 128-4"
$ws.Range("G16").Value = "No. 128-3"
$ws.Range("G17").Value = "No. 128-7"
$ws.Range("G3").Value  = "No. 128-10"
$ws.Range("G7").Value  = "No. 128-5"
$ws.Range("G11").Value = "No. 128-9"

# New column F got a custom width when the product_code column was added
$ws.Columns("F").ColumnWidth = 12.33

# Leave the selection where the author left it after the edit
[void]$ws.Range("G12").Select()
